$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18").Value = "[Euclides-Soldagem, Guilherme-C.L.P., João Bosco-Fundição, Aderci-Fresagem]"
$ws.Range("C18").Value = "[Leandro-M.S.R.A.C., Leonardo-M.Maq.E.I., Leonardo-Retífica, Aderci-CAD/CAM]"
$ws.Range("D18").Value = "[João Bosco-Fundição, Ismail-Metrologia 2, Leandro-M.S.R.A.C., Aderci-CAD/CAM]"
$ws.Range("E18").Value = "Claudinei-Des. Maq. "
$ws.Range("F18").Value = "[Paulo Rob.-Usin. CNC, Guilherme-C.L.P., Guilherme-C. Hidráulica, Leonardo-Retífica]"

$ws.Range("B19").Value = "[Euclides-Soldagem, Guilherme-C.pneumática, Aderci-CAD/CAM, Aderci-Fresagem]"
$ws.Range("D19").Value = "[-, Ismail-Metrologia 2, Leandro-M.S.R.A.C., Guilherme-C.pneumática]"
$ws.Range("E19").Value = "Claudinei-Des. Maq. "
$ws.Range("F19").Value = "[Paulo Rob.-Usin. CNC, Guilherme-C.L.P., Guilherme-C. Hidráulica, Guilherme-C.pneumática]"

$ws.Range("B20").Value = "[Euclides-Soldagem, Guilherme-C.pneumática, Aderci-CAD/CAM, Aderci-Fresagem]"
$ws.Range("C20").Value = "Euclides-Gest. I"
$ws.Range("D20").Value = "[-, Ismail-Metrologia 2, -, -]"
$ws.Range("E20").Value = "Claudinei-Elem"
$ws.Range("F20").Value = "[Paulo Rob.-Usin. CNC, Guilherme-C.L.P., Guilherme-C. Hidráulica, Leonardo-M.Maq.E.I.]"

$ws.Range("B21").Value = "[Euclides-Soldagem, João Bosco-Fundição, Leandro-M.S.R.A.C., Aderci-Fresagem]"
$ws.Range("C21").Value = "Euclides-Gest. I"
$ws.Range("D21").Value = "[João Bosco-Fundição, Ismail-Metrologia 2, Leonardo-Retífica, Leonardo-M.Maq.E.I.]"
$ws.Range("E21").Value = "Claudinei-Elem"
$ws.Range("F21").Value = "[Paulo Rob.-Usin. CNC, Leonardo-Retífica, Guilherme-C. Hidráulica, Leonardo-M.Maq.E.I.]"
